# Translation sheet ("Translation"): the TextArea that used to be labelled
# "bomb cnt" (row 32) is simplified to just show its numeric value, and two
# brand-new rows are appended for the new "unclicked count" TextArea, which
# the commit message says is UI-only for now ("function not implmented").
#
# Row 38 -> SingleUseId38 (a still-unused placeholder TextArea, GB = "<value>")
# Row 39 -> SingleUseId39 (the new unclick-count TextArea, GB defaults to "0")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# F32 held the literal text "bomb cnt"; it now just shows "0".
# Copy from F4 (already a *text* "0") instead of typing "0" directly, so the
# cell keeps its General style and the value stays text instead of being
# auto-coerced to a number by the Value setter.
$ws.Range("F4").Copy($ws.Range("F32"))

# New row 38.
$ws.Range("B38").Value = "SingleUseId38"
$ws.Range("C38").Value = "Default"
$ws.Range("D38").Value = "Right"
$ws.Range("E38").Value = "LTR"
$ws.Range("F38").Value = "<value>"

# New row 39 - the new "unclick cnt" TextArea.
$ws.Range("B39").Value = "SingleUseId39"
$ws.Range("C39").Value = "Default"
$ws.Range("D39").Value = "Left"
$ws.Range("E39").Value = "LTR"
$ws.Range("F4").Copy($ws.Range("F39"))
